$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.136.13'
$ws.Range("E2").Value = '  -7.88%  '
$ws.Range("D3").Value = '2.856.45'
$ws.Range("E3").Value = '  -10.86%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '470.05'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -12.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '125.08'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -7.28%  '
$ws.Range("D8").Value = '2.852.62'
$ws.Range("E8").Value = '  -10.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.400'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -12.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.60'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -11.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0952'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -16.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.328'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -16.71%  '
$ws.Range("E13").Value = '  -4.64%  '
$ws.Range("D14").Value = '3.348.32'
$ws.Range("E14").Value = '  -10.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.15'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -10.36%  '
$ws.Range("D16").Value = '54.173.20'
$ws.Range("E16").Value = '  -7.87%  '
$ws.Range("D17").Value = '2.862.22'
$ws.Range("E17").Value = '  -10.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000133'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -15.09%  '
$ws.Range("E19").Value = '  -10.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.39'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -13.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.03'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -13.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '291.57'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -18.90%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.442'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -14.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '58.56'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -16.23%  '
$ws.Range("E26").Value = '  +0.43%  '
$ws.Range("E27").Value = '  -10.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("D29").Value = '0.0₃0797'
$ws.Range("E29").Value = '  -16.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.17'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -12.83%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.19'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -11.95%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.12'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -6.77%  '
$ws.Range("E33").Value = '  -16.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.74'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -13.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.20'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -14.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '135.12'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -16.24%  '
$ws.Range("E37").Value = '  -14.88%  '
$ws.Range("E38").Value = '  -14.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.88'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -12.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0612'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -12.91%  '
$ws.Range("D41").Value = '2.883.52'
$ws.Range("E41").Value = '  -10.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '34.98'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -14.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.944'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -13.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.598'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -16.17%  '
$ws.Range("E46").Value = '  -12.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.38'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -15.79%  '
$ws.Range("D48").Value = '2.036.84'
$ws.Range("E48").Value = '  -11.35%  '
$ws.Range("E49").Value = '  -15.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.82'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -13.55%  '
$ws.Range("E51").Value = '  -11.87%  '
